# Update market-price-derived columns (H-N) across all item sheets.
# Values below were recomputed by the scheduled pricing runner; this
# script re-applies the refreshed figures (and, where a Leve no longer
# has an HQ/NQ variant, adds/removes the corresponding profit cell).

$wb = $excel.ActiveWorkbook

$changes = @(
    @{ Sheet = "ALC"; Row = 17; Col = 8; Value = 1833.9166 },
    @{ Sheet = "ALC"; Row = 17; Col = 9; Value = 1333.3334 },
    @{ Sheet = "ALC"; Row = 17; Col = 10; Value = 2084.2083 },
    @{ Sheet = "ALC"; Row = 17; Col = 11; Value = 4000.0002 },
    @{ Sheet = "ALC"; Row = 17; Col = 12; Value = 6252.624899999999 },
    @{ Sheet = "ALC"; Row = 17; Col = 13; Value = -3832.0002 },
    @{ Sheet = "ALC"; Row = 17; Col = 14; Value = -6588.624899999999 },
    @{ Sheet = "ALC"; Row = 19; Col = 8; Value = 500.2353 },
    @{ Sheet = "ALC"; Row = 19; Col = 9; Value = 669.25 },
    @{ Sheet = "ALC"; Row = 19; Col = 10; Value = 350 },
    @{ Sheet = "ALC"; Row = 19; Col = 11; Value = 669.25 },
    @{ Sheet = "ALC"; Row = 19; Col = 12; Value = 350 },
    @{ Sheet = "ALC"; Row = 19; Col = 13; Value = -494.25 },
    @{ Sheet = "ALC"; Row = 19; Col = 14; Value = -700 },
    @{ Sheet = "ALC"; Row = 32; Col = 8; Value = 0 },
    @{ Sheet = "ALC"; Row = 32; Col = 10; Value = 0 },
    @{ Sheet = "ALC"; Row = 32; Col = 12; Value = 0 },
    @{ Sheet = "ALC"; Row = 32; Col = 14; Value = $null },
    @{ Sheet = "ALC"; Row = 40; Col = 8; Value = 6110.5 },
    @{ Sheet = "ALC"; Row = 40; Col = 10; Value = 9999 },
    @{ Sheet = "ALC"; Row = 40; Col = 12; Value = 9999 },
    @{ Sheet = "ALC"; Row = 40; Col = 14; Value = -10349 },
    @{ Sheet = "ALC"; Row = 51; Col = 8; Value = 4400 },
    @{ Sheet = "ALC"; Row = 51; Col = 9; Value = 3000 },
    @{ Sheet = "ALC"; Row = 51; Col = 10; Value = 4750 },
    @{ Sheet = "ALC"; Row = 51; Col = 11; Value = 3000 },
    @{ Sheet = "ALC"; Row = 51; Col = 12; Value = 4750 },
    @{ Sheet = "ALC"; Row = 51; Col = 13; Value = -2516 },
    @{ Sheet = "ALC"; Row = 51; Col = 14; Value = -5718 },
    @{ Sheet = "ALC"; Row = 74; Col = 8; Value = 11391.5 },
    @{ Sheet = "ALC"; Row = 74; Col = 9; Value = 4587.25 },
    @{ Sheet = "ALC"; Row = 74; Col = 11; Value = 4587.25 },
    @{ Sheet = "ALC"; Row = 74; Col = 13; Value = -3651.25 },
    @{ Sheet = "ALC"; Row = 77; Col = 8; Value = 11391.5 },
    @{ Sheet = "ALC"; Row = 77; Col = 9; Value = 4587.25 },
    @{ Sheet = "ALC"; Row = 77; Col = 11; Value = 22936.25 },
    @{ Sheet = "ALC"; Row = 77; Col = 13; Value = -18256.25 },
    @{ Sheet = "ALC"; Row = 98; Col = 8; Value = 1403.4667 },
    @{ Sheet = "ALC"; Row = 98; Col = 9; Value = 1403.4667 },
    @{ Sheet = "ALC"; Row = 98; Col = 10; Value = 0 },
    @{ Sheet = "ALC"; Row = 98; Col = 11; Value = 1403.4667 },
    @{ Sheet = "ALC"; Row = 98; Col = 12; Value = 0 },
    @{ Sheet = "ALC"; Row = 98; Col = 13; Value = 94.53330000000005 },
    @{ Sheet = "ALC"; Row = 98; Col = 14; Value = $null },
    @{ Sheet = "ALC"; Row = 122; Col = 8; Value = 1403.4667 },
    @{ Sheet = "ALC"; Row = 122; Col = 9; Value = 1403.4667 },
    @{ Sheet = "ALC"; Row = 122; Col = 10; Value = 0 },
    @{ Sheet = "ALC"; Row = 122; Col = 11; Value = 4210.4001 },
    @{ Sheet = "ALC"; Row = 122; Col = 12; Value = 0 },
    @{ Sheet = "ALC"; Row = 122; Col = 13; Value = -1760.4001 },
    @{ Sheet = "ALC"; Row = 122; Col = 14; Value = $null },
    @{ Sheet = "ALC"; Row = 131; Col = 8; Value = 1084 },
    @{ Sheet = "ALC"; Row = 131; Col = 9; Value = 591.3333 },
    @{ Sheet = "ALC"; Row = 131; Col = 11; Value = 1773.9999 },
    @{ Sheet = "ALC"; Row = 131; Col = 13; Value = 3266.0001 },
    @{ Sheet = "ALC"; Row = 132; Col = 8; Value = 13374.353 },
    @{ Sheet = "ALC"; Row = 132; Col = 9; Value = 16820.691 },
    @{ Sheet = "ALC"; Row = 132; Col = 10; Value = 2173.75 },
    @{ Sheet = "ALC"; Row = 132; Col = 11; Value = 50462.073 },
    @{ Sheet = "ALC"; Row = 132; Col = 12; Value = 6521.25 },
    @{ Sheet = "ALC"; Row = 132; Col = 13; Value = -47932.073 },
    @{ Sheet = "ALC"; Row = 132; Col = 14; Value = -11581.25 },
    @{ Sheet = "ALC"; Row = 138; Col = 8; Value = 1793.9524 },
    @{ Sheet = "ALC"; Row = 138; Col = 9; Value = 1511.2 },
    @{ Sheet = "ALC"; Row = 138; Col = 10; Value = 2500.8333 },
    @{ Sheet = "ALC"; Row = 138; Col = 11; Value = 4533.6 },
    @{ Sheet = "ALC"; Row = 138; Col = 12; Value = 7502.499899999999 },
    @{ Sheet = "ALC"; Row = 138; Col = 13; Value = 606.3999999999996 },
    @{ Sheet = "ALC"; Row = 138; Col = 14; Value = -17782.4999 },
    @{ Sheet = "ARM"; Row = 61; Col = 8; Value = 8500 },
    @{ Sheet = "ARM"; Row = 61; Col = 9; Value = 8500 },
    @{ Sheet = "ARM"; Row = 61; Col = 11; Value = 8500 },
    @{ Sheet = "ARM"; Row = 61; Col = 13; Value = -8288 },
    @{ Sheet = "ARM"; Row = 74; Col = 8; Value = 4554.227 },
    @{ Sheet = "ARM"; Row = 74; Col = 9; Value = 3185.2666 },
    @{ Sheet = "ARM"; Row = 74; Col = 10; Value = 7487.7144 },
    @{ Sheet = "ARM"; Row = 74; Col = 11; Value = 3185.2666 },
    @{ Sheet = "ARM"; Row = 74; Col = 12; Value = 7487.7144 },
    @{ Sheet = "ARM"; Row = 74; Col = 13; Value = -2311.2666 },
    @{ Sheet = "ARM"; Row = 74; Col = 14; Value = -9235.714400000001 },
    @{ Sheet = "ARM"; Row = 77; Col = 8; Value = 4554.227 },
    @{ Sheet = "ARM"; Row = 77; Col = 9; Value = 3185.2666 },
    @{ Sheet = "ARM"; Row = 77; Col = 10; Value = 7487.7144 },
    @{ Sheet = "ARM"; Row = 77; Col = 11; Value = 15926.333 },
    @{ Sheet = "ARM"; Row = 77; Col = 12; Value = 37438.572 },
    @{ Sheet = "ARM"; Row = 77; Col = 13; Value = -11558.333 },
    @{ Sheet = "ARM"; Row = 77; Col = 14; Value = -46174.572 },
    @{ Sheet = "ARM"; Row = 88; Col = 8; Value = 2369.8 },
    @{ Sheet = "ARM"; Row = 88; Col = 10; Value = 2462.25 },
    @{ Sheet = "ARM"; Row = 88; Col = 12; Value = 2462.25 },
    @{ Sheet = "ARM"; Row = 88; Col = 14; Value = -3274.25 },
    @{ Sheet = "ARM"; Row = 91; Col = 8; Value = 2369.8 },
    @{ Sheet = "ARM"; Row = 91; Col = 10; Value = 2462.25 },
    @{ Sheet = "ARM"; Row = 91; Col = 12; Value = 2462.25 },
    @{ Sheet = "ARM"; Row = 91; Col = 14; Value = -5270.25 },
    @{ Sheet = "ARM"; Row = 92; Col = 8; Value = 50550 },
    @{ Sheet = "ARM"; Row = 92; Col = 10; Value = 50550 },
    @{ Sheet = "ARM"; Row = 92; Col = 12; Value = 50550 },
    @{ Sheet = "ARM"; Row = 92; Col = 14; Value = -55542 },
    @{ Sheet = "ARM"; Row = 110; Col = 8; Value = 4413.75 },
    @{ Sheet = "ARM"; Row = 110; Col = 9; Value = 3968.3333 },
    @{ Sheet = "ARM"; Row = 110; Col = 11; Value = 3968.3333 },
    @{ Sheet = "ARM"; Row = 110; Col = 13; Value = -1923.3333 },
    @{ Sheet = "ARM"; Row = 122; Col = 8; Value = 2115.5 },
    @{ Sheet = "ARM"; Row = 122; Col = 9; Value = 1898 },
    @{ Sheet = "ARM"; Row = 122; Col = 11; Value = 5694 },
    @{ Sheet = "ARM"; Row = 122; Col = 13; Value = -3244 },
    @{ Sheet = "ARM"; Row = 132; Col = 8; Value = 3110.6667 },
    @{ Sheet = "ARM"; Row = 132; Col = 9; Value = 1856.7142 },
    @{ Sheet = "ARM"; Row = 132; Col = 10; Value = 7499.5 },
    @{ Sheet = "ARM"; Row = 132; Col = 11; Value = 5570.142599999999 },
    @{ Sheet = "ARM"; Row = 132; Col = 12; Value = 22498.5 },
    @{ Sheet = "ARM"; Row = 132; Col = 13; Value = -3040.142599999999 },
    @{ Sheet = "ARM"; Row = 132; Col = 14; Value = -27558.5 },
    @{ Sheet = "ARM"; Row = 136; Col = 8; Value = 8500 },
    @{ Sheet = "ARM"; Row = 136; Col = 9; Value = 8500 },
    @{ Sheet = "ARM"; Row = 136; Col = 11; Value = 25500 },
    @{ Sheet = "ARM"; Row = 136; Col = 13; Value = -22950 },
    @{ Sheet = "BSM"; Row = 86; Col = 8; Value = 8685.429 },
    @{ Sheet = "BSM"; Row = 86; Col = 10; Value = 9999.666999999999 },
    @{ Sheet = "BSM"; Row = 86; Col = 12; Value = 9999.666999999999 },
    @{ Sheet = "BSM"; Row = 86; Col = 14; Value = -12245.667 },
    @{ Sheet = "BSM"; Row = 89; Col = 8; Value = 8685.429 },
    @{ Sheet = "BSM"; Row = 89; Col = 10; Value = 9999.666999999999 },
    @{ Sheet = "BSM"; Row = 89; Col = 12; Value = 49998.335 },
    @{ Sheet = "BSM"; Row = 89; Col = 14; Value = -61230.335 },
    @{ Sheet = "BSM"; Row = 134; Col = 8; Value = 3444.4285 },
    @{ Sheet = "BSM"; Row = 134; Col = 9; Value = 2977.0908 },
    @{ Sheet = "BSM"; Row = 134; Col = 11; Value = 8931.2724 },
    @{ Sheet = "BSM"; Row = 134; Col = 13; Value = -6396.2724 },
    @{ Sheet = "CRP"; Row = 12; Col = 8; Value = 840 },
    @{ Sheet = "CRP"; Row = 12; Col = 9; Value = 409.16666 },
    @{ Sheet = "CRP"; Row = 12; Col = 10; Value = 1486.25 },
    @{ Sheet = "CRP"; Row = 12; Col = 11; Value = 409.16666 },
    @{ Sheet = "CRP"; Row = 12; Col = 12; Value = 1486.25 },
    @{ Sheet = "CRP"; Row = 12; Col = 13; Value = -239.16666 },
    @{ Sheet = "CRP"; Row = 12; Col = 14; Value = -1826.25 },
    @{ Sheet = "CRP"; Row = 31; Col = 8; Value = 6867.5 },
    @{ Sheet = "CRP"; Row = 31; Col = 9; Value = 1663.625 },
    @{ Sheet = "CRP"; Row = 31; Col = 10; Value = 9841.143 },
    @{ Sheet = "CRP"; Row = 31; Col = 11; Value = 1663.625 },
    @{ Sheet = "CRP"; Row = 31; Col = 12; Value = 9841.143 },
    @{ Sheet = "CRP"; Row = 31; Col = 13; Value = -1368.625 },
    @{ Sheet = "CRP"; Row = 31; Col = 14; Value = -10431.143 },
    @{ Sheet = "CRP"; Row = 34; Col = 8; Value = 6867.5 },
    @{ Sheet = "CRP"; Row = 34; Col = 9; Value = 1663.625 },
    @{ Sheet = "CRP"; Row = 34; Col = 10; Value = 9841.143 },
    @{ Sheet = "CRP"; Row = 34; Col = 11; Value = 1663.625 },
    @{ Sheet = "CRP"; Row = 34; Col = 12; Value = 9841.143 },
    @{ Sheet = "CRP"; Row = 34; Col = 13; Value = -1461.625 },
    @{ Sheet = "CRP"; Row = 34; Col = 14; Value = -10245.143 },
    @{ Sheet = "CRP"; Row = 58; Col = 8; Value = 2569.0667 },
    @{ Sheet = "CRP"; Row = 58; Col = 10; Value = 3676.7144 },
    @{ Sheet = "CRP"; Row = 58; Col = 12; Value = 3676.7144 },
    @{ Sheet = "CRP"; Row = 58; Col = 14; Value = -4082.7144 },
    @{ Sheet = "CRP"; Row = 92; Col = 8; Value = 36149.75 },
    @{ Sheet = "CRP"; Row = 92; Col = 10; Value = 36149.75 },
    @{ Sheet = "CRP"; Row = 92; Col = 12; Value = 36149.75 },
    @{ Sheet = "CRP"; Row = 92; Col = 14; Value = -41141.75 },
    @{ Sheet = "CRP"; Row = 132; Col = 8; Value = 4213.4 },
    @{ Sheet = "CRP"; Row = 132; Col = 9; Value = 3611.9333 },
    @{ Sheet = "CRP"; Row = 132; Col = 10; Value = 6017.8 },
    @{ Sheet = "CRP"; Row = 132; Col = 11; Value = 10835.7999 },
    @{ Sheet = "CRP"; Row = 132; Col = 12; Value = 18053.4 },
    @{ Sheet = "CRP"; Row = 132; Col = 13; Value = -8305.7999 },
    @{ Sheet = "CRP"; Row = 132; Col = 14; Value = -23113.4 },
    @{ Sheet = "CRP"; Row = 136; Col = 8; Value = 2569.0667 },
    @{ Sheet = "CRP"; Row = 136; Col = 10; Value = 3676.7144 },
    @{ Sheet = "CRP"; Row = 136; Col = 12; Value = 11030.1432 },
    @{ Sheet = "CRP"; Row = 136; Col = 14; Value = -16130.1432 },
    @{ Sheet = "CUL"; Row = 14; Col = 8; Value = 195.5 },
    @{ Sheet = "CUL"; Row = 14; Col = 9; Value = 195.5 },
    @{ Sheet = "CUL"; Row = 14; Col = 11; Value = 586.5 },
    @{ Sheet = "CUL"; Row = 14; Col = 13; Value = -413.5 },
    @{ Sheet = "CUL"; Row = 59; Col = 8; Value = 1500 },
    @{ Sheet = "CUL"; Row = 59; Col = 9; Value = 0 },
    @{ Sheet = "CUL"; Row = 59; Col = 10; Value = 1500 },
    @{ Sheet = "CUL"; Row = 59; Col = 11; Value = 0 },
    @{ Sheet = "CUL"; Row = 59; Col = 12; Value = 4500 },
    @{ Sheet = "CUL"; Row = 59; Col = 13; Value = $null },
    @{ Sheet = "CUL"; Row = 59; Col = 14; Value = -5580 },
    @{ Sheet = "CUL"; Row = 134; Col = 8; Value = 7622.1113 },
    @{ Sheet = "CUL"; Row = 134; Col = 9; Value = 1266.6666 },
    @{ Sheet = "CUL"; Row = 134; Col = 11; Value = 3799.9998 },
    @{ Sheet = "CUL"; Row = 134; Col = 13; Value = 1270.0002 },
    @{ Sheet = "CUL"; Row = 137; Col = 8; Value = 1250 },
    @{ Sheet = "CUL"; Row = 137; Col = 10; Value = 1500 },
    @{ Sheet = "CUL"; Row = 137; Col = 12; Value = 4500 },
    @{ Sheet = "CUL"; Row = 137; Col = 14; Value = -14700 },
    @{ Sheet = "GSM"; Row = 95; Col = 8; Value = 22274.166 },
    @{ Sheet = "GSM"; Row = 95; Col = 10; Value = 22274.166 },
    @{ Sheet = "GSM"; Row = 95; Col = 12; Value = 22274.166 },
    @{ Sheet = "GSM"; Row = 95; Col = 14; Value = -27766.166 },
    @{ Sheet = "GSM"; Row = 122; Col = 8; Value = 164054.39 },
    @{ Sheet = "GSM"; Row = 122; Col = 10; Value = 4116.5 },
    @{ Sheet = "GSM"; Row = 122; Col = 12; Value = 12349.5 },
    @{ Sheet = "GSM"; Row = 122; Col = 14; Value = -17249.5 },
    @{ Sheet = "GSM"; Row = 126; Col = 8; Value = 3698.1333 },
    @{ Sheet = "GSM"; Row = 126; Col = 9; Value = 3472.6667 },
    @{ Sheet = "GSM"; Row = 126; Col = 11; Value = 10418.0001 },
    @{ Sheet = "GSM"; Row = 126; Col = 13; Value = -7948.000100000001 },
    @{ Sheet = "GSM"; Row = 132; Col = 8; Value = 14691.066 },
    @{ Sheet = "GSM"; Row = 132; Col = 10; Value = 10836 },
    @{ Sheet = "GSM"; Row = 132; Col = 12; Value = 32508 },
    @{ Sheet = "GSM"; Row = 132; Col = 14; Value = -37568 },
    @{ Sheet = "LTW"; Row = 82; Col = 8; Value = 3024.4119 },
    @{ Sheet = "LTW"; Row = 82; Col = 9; Value = 2099.1428 },
    @{ Sheet = "LTW"; Row = 82; Col = 10; Value = 3672.1 },
    @{ Sheet = "LTW"; Row = 82; Col = 11; Value = 2099.1428 },
    @{ Sheet = "LTW"; Row = 82; Col = 12; Value = 3672.1 },
    @{ Sheet = "LTW"; Row = 82; Col = 13; Value = -1738.1428 },
    @{ Sheet = "LTW"; Row = 82; Col = 14; Value = -4394.1 },
    @{ Sheet = "LTW"; Row = 85; Col = 8; Value = 3024.4119 },
    @{ Sheet = "LTW"; Row = 85; Col = 9; Value = 2099.1428 },
    @{ Sheet = "LTW"; Row = 85; Col = 10; Value = 3672.1 },
    @{ Sheet = "LTW"; Row = 85; Col = 11; Value = 2099.1428 },
    @{ Sheet = "LTW"; Row = 85; Col = 12; Value = 3672.1 },
    @{ Sheet = "LTW"; Row = 85; Col = 13; Value = -851.1428000000001 },
    @{ Sheet = "LTW"; Row = 85; Col = 14; Value = -6168.1 },
    @{ Sheet = "LTW"; Row = 136; Col = 8; Value = 3994.182 },
    @{ Sheet = "LTW"; Row = 136; Col = 9; Value = 3797.7 },
    @{ Sheet = "LTW"; Row = 136; Col = 11; Value = 11393.1 },
    @{ Sheet = "LTW"; Row = 136; Col = 13; Value = -8843.099999999999 },
    @{ Sheet = "WVR"; Row = 39; Col = 8; Value = 50000 },
    @{ Sheet = "WVR"; Row = 39; Col = 10; Value = 50000 },
    @{ Sheet = "WVR"; Row = 39; Col = 12; Value = 50000 },
    @{ Sheet = "WVR"; Row = 39; Col = 14; Value = -50826 },
    @{ Sheet = "WVR"; Row = 107; Col = 8; Value = 635.2222 },
    @{ Sheet = "WVR"; Row = 107; Col = 9; Value = 635.2222 },
    @{ Sheet = "WVR"; Row = 107; Col = 10; Value = 0 },
    @{ Sheet = "WVR"; Row = 107; Col = 11; Value = 1905.6666 },
    @{ Sheet = "WVR"; Row = 107; Col = 12; Value = 0 },
    @{ Sheet = "WVR"; Row = 107; Col = 13; Value = 14.33339999999998 },
    @{ Sheet = "WVR"; Row = 107; Col = 14; Value = $null },
    @{ Sheet = "WVR"; Row = 113; Col = 8; Value = 1394.4 },
    @{ Sheet = "WVR"; Row = 113; Col = 9; Value = 1543 },
    @{ Sheet = "WVR"; Row = 113; Col = 10; Value = 800 },
    @{ Sheet = "WVR"; Row = 113; Col = 11; Value = 4629 },
    @{ Sheet = "WVR"; Row = 113; Col = 12; Value = 2400 },
    @{ Sheet = "WVR"; Row = 113; Col = 13; Value = -2459 },
    @{ Sheet = "WVR"; Row = 113; Col = 14; Value = -6740 },
    @{ Sheet = "WVR"; Row = 132; Col = 8; Value = 2034 },
    @{ Sheet = "WVR"; Row = 132; Col = 9; Value = 1937.2106 },
    @{ Sheet = "WVR"; Row = 132; Col = 11; Value = 5811.6318 },
    @{ Sheet = "WVR"; Row = 132; Col = 13; Value = -3281.6318 },
    @{ Sheet = "WVR"; Row = 136; Col = 8; Value = 2631.7778 },
    @{ Sheet = "WVR"; Row = 136; Col = 9; Value = 1638.8214 },
    @{ Sheet = "WVR"; Row = 136; Col = 13; Value = -2366.4642 }
)

foreach ($c in $changes) {
    $ws = $wb.Worksheets.Item($c.Sheet)
    $ws.Cells.Item($c.Row, $c.Col).Value = $c.Value
}
